$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial for every data row (2..176).
# All of them move from 45221 (2023-10-22) to 45224 (2023-10-25).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 176 }

$ws.Range("C2:C$lastRow").Value = 45224
